$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) cells, forcing text to preserve exact formatting
$priceCells = @{
    "D2" = "67.890.97"
    "D3" = "3.379.30"
    "D5" = "593.05"
    "D6" = "187.58"
    "D8" = "0.597"
    "D13" = "3.921.50"
    "D14" = "639.57"
    "D15" = "8.64"
    "D16" = "67.890.97"
    "D17" = "3.380.87"
    "D18" = "0.119"
    "D19" = "18.09"
    "D20" = "11.15"
    "D21" = "0.911"
    "D22" = "18.00"
    "D24" = "100.05"
    "D27" = "9.84"
    "D28" = "32.75"
    "D30" = "6.95"
    "D31" = "616.19"
    "D33" = "4.047.49"
    "D38" = "2.83"
    "D40" = "33.87"
    "D41" = "3.26"
    "D42" = "0.0₃0705"
    "D44" = "0.345"
    "D48" = "1.39"
    "D49" = "1.00"
    "D50" = "128.11"
    "D51" = "7.74"
}
foreach ($addr in $priceCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $priceCells[$addr]
    $rng.Style = "Normal"
}

# Update other cells (Coin name, Link, Volume%) directly
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("E14").Value = "  +7.87%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  +2.42%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("E26").Value = "  +5.35%  "
$ws.Range("E27").Value = "  +3.71%  "
$ws.Range("E28").Value = "  +6.16%  "
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("E31").Value = "  +5.37%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("E33").Value = "  +7.12%  "
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  +6.59%  "
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("E41").Value = "  +3.02%  "
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("E45").Value = "  +2.21%  "
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E48").Value = "  +11.95%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("E51").Value = "  +4.76%  "
